# New crime data collected — refresh the weekly CompStat figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text (Volume/Number + reporting week dates) ----
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# ---- Row 14 (Murder) ----
$ws.Range("N14").Value = -80

# ---- Row 15 (Rape) ----
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = -7.692307692307
$ws.Range("M15").Value = 9.090909090909
$ws.Range("N15").Value = -68.421052631578

# ---- Row 16 (Robbery) ----
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -27.777777777777
$ws.Range("I16").Value = 132
$ws.Range("J16").Value = 108
$ws.Range("K16").Value = 22.222222222222
$ws.Range("L16").Value = 15.789473684210
$ws.Range("M16").Value = 3.125
$ws.Range("N16").Value = -72.895277207392

# ---- Row 17 (Fel. Assault) ----
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -5.555555555555
$ws.Range("I17").Value = 260
$ws.Range("J17").Value = 225
$ws.Range("K17").Value = 15.555555555555
$ws.Range("L17").Value = 18.721461187214
$ws.Range("M17").Value = 134.234234234234
$ws.Range("N17").Value = -25.072046109510

# ---- Row 18 (Burglary) ----
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 67
$ws.Range("K18").Value = 40.298507462686
$ws.Range("L18").Value = -7.843137254901
$ws.Range("M18").Value = 147.368421052632
$ws.Range("N18").Value = -61.475409836065

# ---- Row 19 (Gr. Larceny) ----
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -6.896551724137
$ws.Range("I19").Value = 214
$ws.Range("J19").Value = 236
$ws.Range("K19").Value = -9.322033898305
$ws.Range("L19").Value = -20.149253731343
$ws.Range("M19").Value = 51.773049645390
$ws.Range("N19").Value = -38.681948424068

# ---- Row 20 (G.L.A.) ----
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 47
$ws.Range("K20").Value = -10.638297872340
$ws.Range("L20").Value = -14.285714285714
$ws.Range("M20").Value = 68
$ws.Range("N20").Value = -81.415929203539

# ---- Row 21 (TOTAL) ----
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -8.080808080808
$ws.Range("I21").Value = 758
$ws.Range("J21").Value = 699
$ws.Range("K21").Value = 8.440629470672
$ws.Range("L21").Value = -1.302083333333
$ws.Range("M21").Value = 66.593406593406
$ws.Range("N21").Value = -55.698421975453

# ---- Row 23 (Petit Larceny) ----
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = -46.153846153846
$ws.Range("F23").Value = 26
$ws.Range("G23").Value = 32
$ws.Range("H23").Value = -18.75
$ws.Range("I23").Value = 224
$ws.Range("J23").Value = 218
$ws.Range("K23").Value = 2.752293577981
$ws.Range("L23").Value = -4.680851063829
$ws.Range("M23").Value = 68.421052631578

# ---- Row 24 (Retail Theft) ----
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -40.909090909090
$ws.Range("F24").Value = 53
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -39.080459770114
$ws.Range("I24").Value = 439
$ws.Range("J24").Value = 524
$ws.Range("K24").Value = -16.221374045801
$ws.Range("L24").Value = -6.196581196581
$ws.Range("M24").Value = 21.270718232044

# ---- Row 25 (Misd. Assault) ----
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -80
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -71.875
$ws.Range("I25").Value = 89
$ws.Range("J25").Value = 158
$ws.Range("K25").Value = -43.670886075949
$ws.Range("L25").Value = -22.608695652173

# ---- Row 26 (UCR Rape*) ----
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -45.454545454545
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = 6.666666666666
$ws.Range("I26").Value = 443
$ws.Range("J26").Value = 339
$ws.Range("K26").Value = 30.678466076696
$ws.Range("L26").Value = 34.650455927051
$ws.Range("M26").Value = 11.868686868686

# ---- Row 27 (Other Sex Crimes) — C/D/E flip to the "no data" placeholders ----
$ws.Range("C27").Formula = "'0"
$ws.Range("A14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Formula = "'0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Formula = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("I27").Value = 18
$ws.Range("K27").Value = -5.263157894736
$ws.Range("L27").Value = 5.882352941176

# ---- Row 28 (Shooting Vic.) ----
$ws.Range("D28").Value = 1
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = 3.030303030303
$ws.Range("L28").Value = 0

# ---- Row 29 (Shooting Inc.) — D/E flip from placeholders to real numbers ----
$ws.Range("D29").Value = 1
$ws.Range("G29").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = -100
$ws.Range("H29").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 12
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -64.705882352941
$ws.Range("N29").Value = -90.769230769230

# ---- Row 30 (Hate Crimes) — D/E flip from placeholders to real numbers ----
$ws.Range("D30").Value = 1
$ws.Range("G30").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = -100
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 11
$ws.Range("K30").Value = -63.636363636363
$ws.Range("L30").Value = -63.636363636363
$ws.Range("N30").Value = -93.220338983050

Write-Output "edit applied"
